$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast-error statistics following the quarterly indexing bug-fix.
# Each row corresponds to quarter Q0..Q9 (rows 2-11); columns B:F are the
# recalculated error statistics and column G is the corrected sample count N.

$data = @{
    2  = @(0.2735553655122697, 0.705428612614435,  2.297651370956758, 1.515800570971247, 1.508984402319302, 42)
    3  = @(0.2729092050941885, 1.324535644936742,  6.788792271302114, 2.605531091985301, 2.623389154891797, 41)
    4  = @(0.313763359840497,  1.409704594726758,  6.400624557163932, 2.529945564071277, 2.542394815147271, 40)
    5  = @(0.3440745136685818, 1.359304277531268,  5.991552565471967, 2.447764810081223, 2.455142007288035, 39)
    6  = @(0.3028255374049251, 1.379900412721921,  6.551300714805526, 2.559550881464466, 2.575690416655235, 38)
    7  = @(0.3261224757496547, 1.493742090381438,  6.991343546469828, 2.644114889045071, 2.66011979711121,  37)
    8  = @(0.2710005396799194, 1.468954151564265,  6.346593684193093, 2.519244665409276, 2.540154585288666, 36)
    9  = @(0.2060309553896736, 1.440504172134808,  7.038220111909374, 2.652964400799485, 2.683566533259154, 35)
    10 = @(0.2191995983370183, 1.401612898746454,  6.624809008690754, 2.573870433547647, 2.603085922297586, 34)
    11 = @(0.2232300813796648, 1.4712938842328,    7.299867736119196, 2.70182674058112,  2.7343371580653,   33)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
